$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "Tractor" dependency block mirroring the existing
# Bike/Car/Farzi sections, introducing two new source-code files worth of
# methods/children: TractorPollutionPermit (fetch_tractor, TractorPesticides)
# and TractorPesticides (fetch_pesticides_permit).

$ws.Range("B27").Value = "TractorPollutionPermit"
$ws.Range("C28").Value = "fetch_tractor"
$ws.Range("C29").Value = "TractorPesticides"
$ws.Range("B30").Value = "TractorPesticides"
$ws.Range("C31").Value = "fetch_pesticides_permit"
